$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to Text format temporarily so numeric-looking strings
# (e.g. "67.469.16", "587.88") are stored as text, matching the source data,
# rather than being auto-converted to numbers by Excel.
$dCol = $ws.Range("D2:D51")
$dCol.NumberFormat = "@"

$ws.Range("D2").Value = '67.469.16'
$ws.Range("E2").Value = '  +4.90%  '
$ws.Range("D3").Value = '3.493.63'
$ws.Range("E3").Value = '  +4.91%  '
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").Value = '587.88'
$ws.Range("E5").Value = '  +6.33%  '
$ws.Range("D6").Value = '186.66'
$ws.Range("E6").Value = '  +7.49%  '
$ws.Range("D7").Value = '0.636'
$ws.Range("E7").Value = '  +1.00%  '
$ws.Range("D8").Value = '3.489.43'
$ws.Range("E8").Value = '  +5.03%  '
$ws.Range("E9").Value = '  +0.04%  '
$ws.Range("D10").Value = '0.174'
$ws.Range("E10").Value = '  +2.02%  '
$ws.Range("D11").Value = '0.652'
$ws.Range("E11").Value = '  +3.06%  '
$ws.Range("D12").Value = '56.56'
$ws.Range("E12").Value = '  +5.75%  '
$ws.Range("E13").Value = '  +1.13%  '
$ws.Range("D14").Value = '9.50'
$ws.Range("E14").Value = '  +4.42%  '
$ws.Range("D15").Value = '4.048.58'
$ws.Range("E15").Value = '  +4.80%  '
$ws.Range("D16").Value = '19.00'
$ws.Range("E16").Value = '  +4.70%  '
$ws.Range("D17").Value = '3.497.13'
$ws.Range("E17").Value = '  +5.10%  '
$ws.Range("D18").Value = '67.536.85'
$ws.Range("E18").Value = '  +4.82%  '
$ws.Range("D19").Value = '12.24'
$ws.Range("E19").Value = '  +4.06%  '
$ws.Range("E20").Value = '  -0.91%  '
$ws.Range("D22").Value = '488.65'
$ws.Range("E22").Value = '  +7.98%  '
$ws.Range("D23").Value = '5.39'
$ws.Range("E23").Value = '  +6.68%  '
$ws.Range("D24").Value = '16.65'
$ws.Range("E24").Value = '  +19.79%  '
$ws.Range("E25").Value = '  +9.44%  '
$ws.Range("D26").Value = '90.18'
$ws.Range("E26").Value = '  +2.83%  '
$ws.Range("D27").Value = '2.96'
$ws.Range("E27").Value = '  +2.66%  '
$ws.Range("D28").Value = '11.03'
$ws.Range("E28").Value = '  +4.18%  '
$ws.Range("D29").Value = '9.18'
$ws.Range("E29").Value = '  +6.88%  '
$ws.Range("D30").Value = '31.58'
$ws.Range("E30").Value = '  +1.40%  '
$ws.Range("E31").Value = '  +10.75%  '
$ws.Range("E32").Value = '  +3.16%  '
$ws.Range("D33").Value = '64.41'
$ws.Range("E33").Value = '  +4.31%  '
$ws.Range("D34").Value = '597.09'
$ws.Range("E34").Value = '  +5.12%  '
$ws.Range("E35").Value = '  +4.78%  '
$ws.Range("E36").Value = '  +6.96%  '
$ws.Range("E37").Value = '  -0.07%  '
$ws.Range("E38").Value = '  +3.99%  '
$ws.Range("D39").Value = '3.57'
$ws.Range("E39").Value = '  +1.74%  '
$ws.Range("E40").Value = '  +5.79%  '
$ws.Range("D41").Value = '0.0₃0770'
$ws.Range("E41").Value = '  +5.58%  '
$ws.Range("D42").Value = '3.267.73'
$ws.Range("E42").Value = '  +6.57%  '
$ws.Range("D43").Value = '2.94'
$ws.Range("E43").Value = '  +7.05%  '
$ws.Range("E44").Value = '  +3.95%  '
$ws.Range("E45").Value = '  +3.65%  '
$ws.Range("B46").Value = 'dogwifhat'
$ws.Range("C46").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D46").Value = '2.78'
$ws.Range("E46").Value = '  +22.93%  '
$ws.Range("B47").Value = 'ApeXProtocol'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D47").Value = '3.26'
$ws.Range("E47").Value = '  +2.37%  '
$ws.Range("E48").Value = '  +1.15%  '
$ws.Range("B49").Value = 'THORChain'
$ws.Range("C49").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D49").Value = '8.82'
$ws.Range("E49").Value = '  +8.04%  '
$ws.Range("B50").Value = 'LidoDAOToken'
$ws.Range("C50").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D50").Value = '3.27'
$ws.Range("E50").Value = '  +13.60%  '
$ws.Range("D51").Value = '1.00'
$ws.Range("E51").Value = '  +0.15%  '

# Restore original (unformatted) cell style now that text values are set,
# so no stray style index is left attached to the cells.
$dCol.ClearFormats()

